$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Add a new header "description" in the next empty column (M1) after the
# existing headers (A1:L1).
$ws.Cells.Item(1, 13).Value = "description"

# Update the active selection to the newly added cell, matching Excel's
# behavior of moving the selection to the cell just edited.
$ws.Range("M1").Select()
